# Fruta / hortaliza, semanal
# Adds a new week of price data (2022-01-07, serial 44568) for
# "Feria Lagunitas de Puerto Montt - Damasco" and pushes the previous
# history down so it stays intact further down the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: push existing rows 29 & 30 down to rows 32 & 33 -------------
# (old row 30 -> new row 33, exact copy, keeps its own values/style)
$ws.Range("A30:T30").Copy($ws.Range("A33:T33"))
# (old row 29 -> new row 32, exact copy, keeps its own values/style)
$ws.Range("A29:T29").Copy($ws.Range("A32:T32"))

# --- Step 2: seed rows 31 and 30 from row 29's original formatting -------
$ws.Range("A29:T29").Copy($ws.Range("A31:T31"))
$ws.Range("A29:T29").Copy($ws.Range("A30:T30"))

# --- Step 3: write the new week's data ------------------------------------

# Row 29: Especial, 200 @ 21000 (18 kilo box)
$ws.Cells.Item(29, 4).Value2 = 44568
$ws.Cells.Item(29, 12).Value2 = "Especial"
$ws.Cells.Item(29, 13).Value2 = 200
$ws.Cells.Item(29, 14).Value2 = 21000
$ws.Cells.Item(29, 15).Value2 = 21000
$ws.Cells.Item(29, 16).Value2 = 21000
$ws.Cells.Item(29, 17).Value2 = "$/caja 18 kilos"
$ws.Cells.Item(29, 19).Value2 = 1167
$ws.Cells.Item(29, 20).Value2 = 18

# Row 30: Primera, 200 @ 18000 (18 kilo box)
$ws.Cells.Item(30, 4).Value2 = 44568
$ws.Cells.Item(30, 12).Value2 = "Primera"
$ws.Cells.Item(30, 13).Value2 = 200
$ws.Cells.Item(30, 14).Value2 = 18000
$ws.Cells.Item(30, 15).Value2 = 18000
$ws.Cells.Item(30, 16).Value2 = 18000
$ws.Cells.Item(30, 17).Value2 = "$/caja 18 kilos"
$ws.Cells.Item(30, 19).Value2 = 1000
$ws.Cells.Item(30, 20).Value2 = 18

# Row 31: Segunda, 200 @ 16000 (18 kilo box)
$ws.Cells.Item(31, 4).Value2 = 44568
$ws.Cells.Item(31, 12).Value2 = "Segunda"
$ws.Cells.Item(31, 13).Value2 = 200
$ws.Cells.Item(31, 14).Value2 = 16000
$ws.Cells.Item(31, 15).Value2 = 16000
$ws.Cells.Item(31, 16).Value2 = 16000
$ws.Cells.Item(31, 17).Value2 = "$/caja 18 kilos"
$ws.Cells.Item(31, 19).Value2 = 889
$ws.Cells.Item(31, 20).Value2 = 18

# Rows 32 and 33 already hold the untouched historical rows (29 & 30
# originals), so no further edits are needed there.
